$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (GET /api/articles/:article_id/comments): mark the 422 column as
# handled ("done"), matching the fill style already used for "done" cells
# elsewhere in that column (e.g. I6).
$ws.Range("I6").Copy()
$ws.Range("I12").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I12").Value = "done"

# Row 13 (POST /api/articles/:article_id/comments): newly implemented
# functionality - fill in the Normal/Queries/422 checklist columns.
$ws.Range("C12").Copy()
$ws.Range("C13").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C13").Value = "done"

$ws.Range("D7").Copy()
$ws.Range("D13").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D13").Value = "N/A"

$ws.Range("I6").Copy()
$ws.Range("I13").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I13").Value = "done"

$excel.CutCopyMode = 0

# Update active selection to D13 to match the author's last-edited cell
$ws.Range("D13").Select()
